{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].insertText(\"ini dokumen uji coba coy!!!!\", Word.InsertLocation.replace);\nfor (let i = 1; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$d.Paragraphs(1).Range.Text = \"ini dokumen uji coba coy!!!!\"\n\nfor ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {\n    $d.Paragraphs($i).Range.Delete()\n}\n"}
